# TC004/TC003 data-driven test sheet update.
# Sheet1 moves from a 2-column (username/password) login sheet to an
# 11-column "set customer" record (login fields + new customer fields).
# Sheet2/Sheet3 are not touched directly - their shared-string indices
# shift automatically once the unused "username"/"password" strings are
# replaced below.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- Row 1: header labels -------------------------------------------------
# (write A1/B1 first so the old "username"/"password" shared strings are
# dropped and "userName"/"passWord" land where the diff expects them)
$ws1.Range("A1").Value2 = "userName"
$ws1.Range("B1").Value2 = "passWord"
$ws1.Range("C1").Value2 = "setCustName"
$ws1.Range("D1").Value2 = "setCustGender"
$ws1.Range("E1").Value2 = "setCustDob"
$ws1.Range("F1").Value2 = "setCustAdd"
$ws1.Range("G1").Value2 = "setCustCity"
$ws1.Range("H1").Value2 = "setCustState"
$ws1.Range("I1").Value2 = "setCustPin"
$ws1.Range("J1").Value2 = "setCustMob"
$ws1.Range("K1").Value2 = "setCustPass"

# give the new header cells (C1:K1) the same filled/bordered look as A1:B1
$ws1.Range("A1").Copy()
$ws1.Range("C1:K1").PasteSpecial(-4122)

# ---- Row 2: sample data row ------------------------------------------------
$ws1.Range("A2").Value2 = "mngr191029"
$ws1.Range("B2").Value2 = "tUdegar"
$ws1.Range("C2").Value2 = "Ashutosh"
$ws1.Range("D2").Value2 = "Male"
$ws1.Range("G2").Value2 = "Gurgaon"
$ws1.Range("F2").Value2 = "Vatika G21"
$ws1.Range("H2").Value2 = "Haryana"
$ws1.Range("I2").Value2 = 122004
$ws1.Range("J2").Value2 = 9953229953
$ws1.Range("K2").Value2 = 123456

# DOB as a real date value, formatted like a short date (m/d/yyyy)
$ws1.Range("E2").Value2 = 32826
$ws1.Range("E2").NumberFormat = "m/d/yyyy"

# match the saved selection/active cell
$ws1.Activate() | Out-Null
$ws1.Range("B2").Select() | Out-Null
